$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add new sheet "Test2" right after "Test1" (it becomes the active sheet,
# matching the target workbook's activeTab).
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Test2"

# Test2 mirrors Test1's B3:E6 block, with each value shifted by +6 and a
# single uniform red-font style (no border/fill) instead of the rainbow
# bordered styles used on Test1.
$values = @(
  @(7, 8, 9, 10),
  @(11, 12, 13, 14),
  @(15, 16, 17, 18),
  @(19, 20, 21, 22)
)
for ($r = 0; $r -lt 4; $r++) {
  for ($c = 0; $c -lt 4; $c++) {
    $ws2.Cells.Item(3 + $r, 2 + $c).Value = $values[$r][$c]
  }
}

$rng2 = $ws2.Range("B3:E6")
$rng2.Font.Color = 255

# New defined name mirroring TileMap, pointed at the new sheet.
[void]$wb.Names.Add("TileMap2", "=Test2!`$B`$3:`$E`$6")

# Test1 selection collapses to a single remote cell now that it's no
# longer the active tab.
[void]$ws1.Range("E10").Select()

# Test2 becomes the active sheet/tab with B3:E6 selected (active cell B3).
[void]$ws2.Select()
[void]$rng2.Select()
